# SP-29 BIS-28: add a "Value Unit" metadata row (validator/handler field)
# to the openbis-metadata sheet, and drop the now-redundant "Value Unit"
# row that used to live on the openbis-data sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("openbis-metadata")
$ws2 = $wb.Worksheets.Item("openbis-data")

# ---- openbis-metadata: insert a new "Value Unit" property row ----
# The new row goes right after "Scale" (row 5) and before "Header Format"
# (currently row 6), pushing the remaining rows down by one.
$ws1.Rows.Item(6).Insert()

# Match the formatting used by the other Property/Value/Description/Example
# rows (copy down from the "Scale" row immediately above).
$ws1.Range("A5:D5").Copy()
$ws1.Range("A6:D6").PasteSpecial(-4122)

$ws1.Range("A6").Value = "Value Unit"
$ws1.Range("C6").Value = "One of mM, uM, RatioT1, or RatioCs"
$ws1.Range("D6").Value = "mM"

# The Description cell for the new row uses the same look as its neighbours
# (italic, 14pt, grey Verdana) but was (re)created as its own font entry.
$descFont = $ws1.Range("C6").Font
$descFont.Name = "Verdana"
$descFont.Size = 14
$descFont.Italic = $true
$descFont.Color = 8421505


# Fill in the example/default values for the two properties that previously
# had blank Value cells.
$ws1.Range("B7").Value = "METABOL HYBRID"
$ws1.Range("B9").Value = "C"

# "Start Data Row" moved from 1 to 3 (to make room for the header rows that
# now include Strain / Value Type above the compound table on openbis-data).
$ws1.Range("B8").Value = 3
$ws1.Range("D8").Value = 3

# ---- openbis-data: drop the old "Value Unit" row (now captured on the
# openbis-metadata sheet instead) ----
$ws2.Rows.Item(3).Delete()

# ---- restore selections to match the saved workbook state ----
$ws2.Rows.Item(3).Select()
$ws1.Range("C13").Select()
